# Nudge the positions of several shapes on slide 2 ("Add of workshop
# overview"). Only X or Y (never both, except for the two circular
# arrows) of each shape moves; everything else about the shapes is left
# untouched.
#
# NOTE on precision: Shape.Left/.Top are exposed in points, but the
# underlying store keeps EMU (914400 EMU per inch / 12700 EMU per
# point). A plain "emu / 12700.0" assignment can truncate to one EMU
# below the intended target once it round-trips through the engine's
# float representation, so each value below has been chosen (as the
# closest representable point value) so it reconstructs to the exact
# target EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$sh = $s.Shapes.Item("Rectangle 9")
$sh.Top = 77.81197357177734

$sh = $s.Shapes.Item("Rectangle 11")
$sh.Top = 77.81197357177734

$sh = $s.Shapes.Item("Striped Right Arrow 18")
$sh.Top = 95.95086669921875

$sh = $s.Shapes.Item("Striped Right Arrow 19")
$sh.Top = 95.95086669921875

$sh = $s.Shapes.Item("Cloud 21")
$sh.Left = 180.40567016601562

$sh = $s.Shapes.Item("Can 22")
$sh.Left = 101.03937530517578

$sh = $s.Shapes.Item("Circular Arrow 24")
$sh.Top = 84.67677307128906

$sh = $s.Shapes.Item("Circular Arrow 25")
$sh.Left = 127.86473083496094
$sh.Top = 313.1276550292969

$sh = $s.Shapes.Item("Rectangle 26")
$sh.Top = 306.2631530761719

$sh = $s.Shapes.Item("Striped Right Arrow 27")
$sh.Top = 324.4020690917969

$sh = $s.Shapes.Item("Left Brace 39")
$sh.Left = 319.23089599609375

$sh = $s.Shapes.Item("TextBox 47")
$sh.Left = 341.1703186035156

$sh = $s.Shapes.Item("Down Arrow Callout 49")
$sh.Top = 306.2631530761719
